$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "bad data" row pointing at the vitamin-d-supplement article
# (row 2: about-cancer/coping/feelings/relaxation/vitamin-d-supplement-cancer-prevention).
# Remaining rows shift up to fill the gap.
$ws.Rows.Item(2).Delete()

# Remove the other "bad data" row, which is now the last row (previously row 9:
# about-cancer/coping/feelings/relaxation/hpv-vaccine-presidents-cancer-panel-improving-uptake).
# Remaining rows shift up to fill the gap.
$ws.Rows.Item(8).Delete()

# Append the new replacement data row at the bottom of the table.
$ws.Range("A8").Value = "news-events/cancer-currents-blog/2019/pancreatic-cancer-targeting-kras-indirectly"
$ws.Range("B8").Value = "Blog Post"
$ws.Range("C8").Value = "English"

# Match the updated selection left behind in the sheet.
$ws.Range("A9:A10").Select()
